# Insert a new weekly price record as row 50 on the active sheet,
# pushing the existing rows 50-77 down to 51-78.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(50).Insert()

$ws.Cells.Item(50, 1).Value = 9
$ws.Cells.Item(50, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(50, 3).Value = "Metropolitana"
$ws.Cells.Item(50, 4).Value = 44518
$ws.Cells.Item(50, 5).Value = 13
$ws.Cells.Item(50, 6).Value = 100112022
$ws.Cells.Item(50, 7).Value = "Arveja Verde"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 25
$ws.Cells.Item(50, 11).Value = 15000
$ws.Cells.Item(50, 12).Value = 16000
$ws.Cells.Item(50, 13).Value = 15480
$ws.Cells.Item(50, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(50, 15).Value = "Región Metropolitana"
$ws.Cells.Item(50, 16).Value = 619
$ws.Cells.Item(50, 17).Value = 25
$ws.Cells.Item(50, 18).Value = "Hortaliza"
